$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'schubert-winterreise_42'
$ws.Cells.Item(2, 2).Value = 'schubert-winterreise_57'
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = '[[''G:(3,5,b7,b9)'', ''C:min'', ''E:dim7''], [''G#:maj'', ''D#:(3,5,b7,b9)'', ''G#:maj''], [''F#:dim7'', ''G:maj'', ''C:min'']]'
$ws.Cells.Item(2, 5).Value = '[[''A:(3,5,b7,b9)'', ''D:min'', ''F#:dim7''], [''A#:maj'', ''F:(3,5,b7,b9)'', ''A#:maj''], [''G#:dim7'', ''A:maj'', ''D:min'']]'
$ws.Cells.Item(2, 6).Value = '[(11.76, 15.64), (20.54, 23.16), (16.24, 17.66)]'
$ws.Cells.Item(2, 7).Value = '[(13.82, 17.98), (23.16, 26.04), (18.66, 20.14)]'
$ws.Cells.Item(2, 8).Value = ''
$ws.Cells.Item(2, 9).Value = ''

# Row 3
$ws.Cells.Item(3, 1).Value = 'schubert-winterreise_44'
$ws.Cells.Item(3, 2).Value = 'isophonics_205'
$ws.Cells.Item(3, 3).Value = 0.04738805970149254
$ws.Cells.Item(3, 4).Value = '[[''F'', ''A#:min'', ''A#:min'']]'
$ws.Cells.Item(3, 5).Value = '[[''E'', ''A:min'', ''A:min/b7'']]'
$ws.Cells.Item(3, 6).Value = '[(14.52, 20.94)]'
$ws.Cells.Item(3, 7).Value = '[(14.868843, 21.242721)]'
$ws.Cells.Item(3, 8).Value = ''
$ws.Cells.Item(3, 9).Value = 'spotify:track:389QX9Q1eUOEZ19vtzzI9O'

# Row 4
$ws.Cells.Item(4, 1).Value = 'schubert-winterreise_126'
$ws.Cells.Item(4, 2).Value = 'isophonics_212'
$ws.Cells.Item(4, 3).Value = 0.2375
$ws.Cells.Item(4, 4).Value = '[[''D:maj/F#'', ''G:maj'', ''D:maj'']]'
$ws.Cells.Item(4, 5).Value = '[[''D'', ''G'', ''D'']]'
$ws.Cells.Item(4, 6).Value = '[(58.08, 65.66)]'
$ws.Cells.Item(4, 7).Value = '[(46.93228, 54.037586)]'
$ws.Cells.Item(4, 8).Value = ''
$ws.Cells.Item(4, 9).Value = ''

# Row 5
$ws.Cells.Item(5, 1).Value = 'schubert-winterreise_156'
$ws.Cells.Item(5, 2).Value = 'schubert-winterreise_59'
$ws.Cells.Item(5, 3).Value = 0.421195652173913
$ws.Cells.Item(5, 4).Value = '[[''C#:7'', ''F#:maj'', ''F#:maj/A#'', ''C#:7'', ''F#:maj'']]'
$ws.Cells.Item(5, 5).Value = '[[''B:7/A'', ''E:maj/G#'', ''E:maj/B'', ''B:7'', ''E:maj/G#'']]'
$ws.Cells.Item(5, 6).Value = '[(21.44, 25.4)]'
$ws.Cells.Item(5, 7).Value = '[(237.58, 253.72)]'
$ws.Cells.Item(5, 8).Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Cells.Item(5, 9).Value = ''

# Row 6
$ws.Cells.Item(6, 1).Value = 'schubert-winterreise_132'
$ws.Cells.Item(6, 2).Value = 'jaah_32'
$ws.Cells.Item(6, 3).Value = 0.0485006518904824
$ws.Cells.Item(6, 4).Value = '[[''C:7'', ''F:maj'', ''F:maj'']]'
$ws.Cells.Item(6, 5).Value = '[[''D:7'', ''G'', ''G'']]'
$ws.Cells.Item(6, 6).Value = '[(6.38, 12.5)]'
$ws.Cells.Item(6, 7).Value = '[(15.79, 20.07)]'
$ws.Cells.Item(6, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Cells.Item(6, 9).Value = ''

# Row 7
$ws.Cells.Item(7, 1).Value = 'schubert-winterreise_59'
$ws.Cells.Item(7, 2).Value = 'isophonics_15'
$ws.Cells.Item(7, 3).Value = 0.3833333333333333
$ws.Cells.Item(7, 4).Value = '[[''E:maj/G#'', ''A:maj'', ''E:maj'', ''B:maj'']]'
$ws.Cells.Item(7, 5).Value = '[[''C'', ''F'', ''C'', ''G'']]'
$ws.Cells.Item(7, 6).Value = '[(59.22, 71.12)]'
$ws.Cells.Item(7, 7).Value = '[(70.203514, 92.97068)]'
$ws.Cells.Item(7, 8).Value = ''
$ws.Cells.Item(7, 9).Value = ''

# Row 8
$ws.Cells.Item(8, 1).Value = 'isophonics_150'
$ws.Cells.Item(8, 2).Value = 'isophonics_135'
$ws.Cells.Item(8, 3).Value = 0.1441647597254004
$ws.Cells.Item(8, 4).Value = '[[''C'', ''D'', ''G'']]'
$ws.Cells.Item(8, 5).Value = '[[''A'', ''B'', ''E'']]'
$ws.Cells.Item(8, 6).Value = '[(55.151295, 60.155195)]'
$ws.Cells.Item(8, 7).Value = '[(13.393711, 17.886772)]'
$ws.Cells.Item(8, 8).Value = ''
$ws.Cells.Item(8, 9).Value = ''

# Row 9
$ws.Cells.Item(9, 1).Value = 'schubert-winterreise_213'
$ws.Cells.Item(9, 2).Value = 'schubert-winterreise_175'
$ws.Cells.Item(9, 3).Value = 0.2666666666666667
$ws.Cells.Item(9, 4).Value = '[[''F:7'', ''A#:maj/F'', ''F:7'', ''A#:maj/F'']]'
$ws.Cells.Item(9, 5).Value = '[[''F:7'', ''A#:maj/F'', ''F:7'', ''A#:maj'']]'
$ws.Cells.Item(9, 6).Value = '[(93.26, 102.46)]'
$ws.Cells.Item(9, 7).Value = '[(97.94, 108.3)]'
$ws.Cells.Item(9, 8).Value = 'spotify:track:1yerCi2iQCVkdHG6rdRn7R'
$ws.Cells.Item(9, 9).Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

# Row 10
$ws.Cells.Item(10, 1).Value = 'schubert-winterreise_89'
$ws.Cells.Item(10, 2).Value = 'jaah_55'
$ws.Cells.Item(10, 3).Value = 0.1106719367588933
$ws.Cells.Item(10, 4).Value = '[[''D:7'', ''G:maj'', ''G:maj/B'']]'
$ws.Cells.Item(10, 5).Value = '[[''G:7'', ''C'', ''C'']]'
$ws.Cells.Item(10, 6).Value = '[(20.08, 22.56)]'
$ws.Cells.Item(10, 7).Value = '[(47.25, 51.08)]'
$ws.Cells.Item(10, 8).Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Cells.Item(10, 9).Value = ''

# Row 11
$ws.Cells.Item(11, 1).Value = 'isophonics_216'
$ws.Cells.Item(11, 2).Value = 'isophonics_133'
$ws.Cells.Item(11, 3).Value = 0.1421370967741936
$ws.Cells.Item(11, 4).Value = '[[''B:min'', ''E'', ''A''], [''A'', ''D'', ''A'']]'
$ws.Cells.Item(11, 5).Value = '[[''E:min'', ''A'', ''D''], [''D/3'', ''G'', ''D'']]'
$ws.Cells.Item(11, 6).Value = '[(23.795215, 29.925283), (28.009637, 33.826235)]'
$ws.Cells.Item(11, 7).Value = '[(25.00678, 29.410731), (57.251357, 73.909779)]'
$ws.Cells.Item(11, 8).Value = ''
$ws.Cells.Item(11, 9).Value = ''

# Row 12
$ws.Cells.Item(12, 1).Value = 'schubert-winterreise_74'
$ws.Cells.Item(12, 2).Value = 'schubert-winterreise_2'
$ws.Cells.Item(12, 3).Value = 0.3939393939393939
$ws.Cells.Item(12, 4).Value = '[[''F:maj'', ''C:7'', ''F:maj'', ''C:7'', ''F:maj'']]'
$ws.Cells.Item(12, 5).Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Cells.Item(12, 6).Value = '[(63.2, 73.12)]'
$ws.Cells.Item(12, 7).Value = '[(20.56, 26.4)]'
$ws.Cells.Item(12, 8).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Cells.Item(12, 9).Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 13
$ws.Cells.Item(13, 1).Value = 'isophonics_155'
$ws.Cells.Item(13, 2).Value = 'schubert-winterreise_186'
$ws.Cells.Item(13, 3).Value = 0.2380952380952381
$ws.Cells.Item(13, 4).Value = '[[''Ab'', ''Db/5'', ''Ab'', ''Db/5'', ''Ab'']]'
$ws.Cells.Item(13, 5).Value = '[[''F:maj'', ''A#:maj'', ''F:maj'', ''A#:maj'', ''F:maj'']]'
$ws.Cells.Item(13, 6).Value = '[(261.828, 275.8)]'
$ws.Cells.Item(13, 7).Value = '[(118.44, 122.34)]'
$ws.Cells.Item(13, 8).Value = ''
$ws.Cells.Item(13, 9).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 14
$ws.Cells.Item(14, 1).Value = 'schubert-winterreise_17'
$ws.Cells.Item(14, 2).Value = 'schubert-winterreise_130'
$ws.Cells.Item(14, 3).Value = 0.4875
$ws.Cells.Item(14, 4).Value = '[[''D:maj/A'', ''G:maj'', ''D:maj/A'']]'
$ws.Cells.Item(14, 5).Value = '[[''E:maj/G#'', ''A:maj'', ''E:maj'']]'
$ws.Cells.Item(14, 6).Value = '[(138.02, 142.34)]'
$ws.Cells.Item(14, 7).Value = '[(55.58, 62.32)]'
$ws.Cells.Item(14, 8).Value = ''
$ws.Cells.Item(14, 9).Value = ''

# Row 15
$ws.Cells.Item(15, 1).Value = 'schubert-winterreise_203'
$ws.Cells.Item(15, 2).Value = 'schubert-winterreise_5'
$ws.Cells.Item(15, 3).Value = 0.07728085867620751
$ws.Cells.Item(15, 4).Value = '[[''G:min/A#'', ''A:hdim7/D#'', ''D:7'', ''G:min'']]'
$ws.Cells.Item(15, 5).Value = '[[''A:min'', ''B:hdim7/D'', ''E:7'', ''A:min'']]'
$ws.Cells.Item(15, 6).Value = '[(94.68, 110.76)]'
$ws.Cells.Item(15, 7).Value = '[(10.94, 16.32)]'
$ws.Cells.Item(15, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Cells.Item(15, 9).Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'

# Row 16
$ws.Cells.Item(16, 1).Value = 'schubert-winterreise_63'
$ws.Cells.Item(16, 2).Value = 'schubert-winterreise_5'
$ws.Cells.Item(16, 3).Value = 0.07984496124031007
$ws.Cells.Item(16, 4).Value = '[[''G:7'', ''C:min'', ''B:dim7/C'', ''C:min''], [''G:7'', ''C:min'', ''C:min'', ''D:hdim7/C''], [''C:min'', ''B:dim7/C'', ''C:min'', ''B:dim7/C'']]'
$ws.Cells.Item(16, 5).Value = '[[''E:7'', ''A:min'', ''G#:dim7'', ''A:min''], [''E:7'', ''A:min'', ''A:min'', ''B:hdim7/D''], [''A:min'', ''G#:dim7'', ''A:min'', ''G#:dim7'']]'
$ws.Cells.Item(16, 6).Value = '[(78.74, 89.26), (26.48, 34.32), (0.24, 9.6)]'
$ws.Cells.Item(16, 7).Value = '[(19.28, 24.94), (13.6, 19.28), (20.66, 26.4)]'
$ws.Cells.Item(16, 8).Value = ''
$ws.Cells.Item(16, 9).Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'

# Row 17
$ws.Cells.Item(17, 1).Value = 'schubert-winterreise_186'
$ws.Cells.Item(17, 2).Value = 'isophonics_53'
$ws.Cells.Item(17, 3).Value = 0.3142857142857143
$ws.Cells.Item(17, 4).Value = '[[''F:maj'', ''C:7'', ''F:maj''], [''F:maj'', ''A#:maj'', ''F:maj'']]'
$ws.Cells.Item(17, 5).Value = '[[''A/3'', ''E:7'', ''A''], [''A'', ''D'', ''A'']]'
$ws.Cells.Item(17, 6).Value = '[(43.34, 59.34), (118.44, 121.68)]'
$ws.Cells.Item(17, 7).Value = '[(58.557, 62.834), (40.49, 47.86)]'
$ws.Cells.Item(17, 8).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Cells.Item(17, 9).Value = ''
